# Update the "time" (cronograma) worksheet with the latest progress values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time")

# Bring this sheet to the front / make it the active sheet, as it already is.
$ws.Activate()

# Progress (% Avance) updates for several tasks in the schedule.
$ws.Range("F6").Value = 0.9     # Inicio de sesión
$ws.Range("F12").Value = 1      # Reporta 2: Mostrar historial de tratamiento...
$ws.Range("F14").Value = 0.2    # Prueba final de desarrollo
$ws.Range("F15").Value = 0.5    # Corregir bug o actualizar mejoras

# Move the current selection to reflect where the user left off editing.
$ws.Range("F15").Select()

# Recalculate so the cached formula result (F19, the overall % average) updates.
$wb.Application.Calculate()
